# regen sval data to filter save games
# Update the pre-computed per-game stat columns (TB, d2S, K, IP) and the
# resulting "sum" column (G) for each data row. Column F (Win) is left
# untouched since it did not change in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => B, C, D, E, G (new values)
$newValues = @{
    2  = @(3.272327238179451,   1.626987699542094,  0.1496068669990043,  0.5333859586016987,  5.582307763322248)
    3  = @(0.6545652718822623,  1.626987699542094,  3.223369029078222,   0.5333859586016987,  6.038307959104277)
    4  = @(1.445647641019636,   1.626987699542094,  0.1496068669990043,  0.5333859586016987,  3.755628166162433)
    5  = @(1.445647641019636,   1.626987699542094,  3.223369029078222,   0.5333859586016987,  6.82939032824165)
    6  = @(1.445647641019636,   1.626987699542094,  0.1496068669990043,  0.5333859586016987,  3.755628166162433)
    7  = @(0.2881169905109251,  0.3048912486333797, 0.7210945179870265,  0.5333859586016987,  1.84748871573303)
    8  = @(3.272327238179451,   1.626987699542094,  3.223369029078222,   0.5333859586016987,  8.656069925401464)
    9  = @(0.6545652718822623,  1.626987699542094,  0.7210945179870265,  0.5333859586016987,  3.536033448013082)
    10 = @(0.04172184405617529, 0.3048912486333797, 0.7210945179870265,  13.86384647080068,   14.93155408147727)
    11 = @(3.272327238179451,   1.626987699542094,  0.7210945179870265,  0.5333859586016987,  6.15379541431027)
    12 = @(3.272327238179451,   1.626987699542094,  0.7210945179870265,  0.5333859586016987,  6.15379541431027)
    13 = @(0.2881169905109251,  0.3048912486333797, 18.71679738969934,   0.5333859586016987,  19.84319158744534)
    14 = @(1.445647641019636,   1.626987699542094,  0.7210945179870265,  0.5333859586016987,  4.327115817150455)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
